# Update Name of Algo
# Applies updated KNN-imputed values for columns C and D on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -7.509
$ws.Range("C12").Value = -10.89
$ws.Range("D14").Value = -7.686000000000002
$ws.Range("D26").Value = -8.026999999999999
$ws.Range("D31").Value = -8.218999999999999
$ws.Range("C32").Value = -13.584
$ws.Range("D35").Value = -7.935
$ws.Range("C36").Value = -12.732
$ws.Range("D37").Value = -7.712999999999999
$ws.Range("C38").Value = -12.703
$ws.Range("D45").Value = -7.603
$ws.Range("C46").Value = -14.162
$ws.Range("C54").Value = -12.705
$ws.Range("C55").Value = -13.624
$ws.Range("D57").Value = -8.074000000000002
$ws.Range("C67").Value = -11.705
$ws.Range("C69").Value = -11.038
$ws.Range("C72").Value = -11.555
$ws.Range("C91").Value = -10.879
$ws.Range("C99").Value = -12.635
$ws.Range("D100").Value = -8.272
$ws.Range("D102").Value = -7.790000000000001
